$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing test case row (Test Case No 4 -> 1, employees=5 -> employees=2) ---
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "employees=2"

# --- Add new test case rows (row 16 = Test Case 2, row 18 = Test Case 3) ---
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "to get help"
$ws.Range("C16").Value = "argv[0] -h"
$ws.Range("D16").Value = "enter inputs"
$ws.Range("E16").Value = "enter inputs"
$ws.Range("F16").Value = "PASS"

$ws.Range("A18").Value = 3
$ws.Range("B18").Value = "with 0 as input or negative input"
$ws.Range("C18").Value = "employees=0 or employee=-10"
$ws.Range("D18").Value = "ID           Name        Type"
$ws.Range("E18").Value = "ID           Name        Type"
$ws.Range("F18").Value = "PASS"

# --- Column C widened to fit the longer text that was just entered ---
$null = $ws.Columns.Item(3).AutoFit()

# --- Selection moved to D20 ---
$null = $ws.Range("D20").Select()
